$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = "27.162.15"
$ws.Cells.Item(2, 5).Value = "  +0.79%  "
$ws.Cells.Item(3, 4).Value = "1.681.11"
$ws.Cells.Item(3, 5).Value = "  +0.40%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(5, 4) "215.36"
$ws.Cells.Item(5, 5).Value = "  +0.30%  "
$ws.Cells.Item(6, 5).Value = "  +0.24%  "
$ws.Cells.Item(7, 5).Value = "  -0.06%  "
$ws.Cells.Item(8, 5).Value = "  +2.17%  "
Set-TextValue $ws.Cells.Item(9, 4) "21.51"
$ws.Cells.Item(9, 5).Value = "  +5.75%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.0624"
$ws.Cells.Item(11, 5).Value = "  +0.35%  "
$ws.Cells.Item(12, 4).Value = "1.916.60"
$ws.Cells.Item(12, 5).Value = "  +0.34%  "
$ws.Cells.Item(13, 4).Value = "1.675.57"
$ws.Cells.Item(13, 5).Value = "  -0.15%  "
Set-TextValue $ws.Cells.Item(14, 4) "4.15"
$ws.Cells.Item(14, 5).Value = "  +1.57%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.537"
$ws.Cells.Item(15, 5).Value = "  +1.94%  "
Set-TextValue $ws.Cells.Item(16, 4) "66.39"
$ws.Cells.Item(17, 4).Value = "27.140.51"
$ws.Cells.Item(17, 5).Value = "  +0.63%  "
Set-TextValue $ws.Cells.Item(18, 4) "238.68"
$ws.Cells.Item(18, 5).Value = "  +0.61%  "
$ws.Cells.Item(19, 5).Value = "  +0.45%  "
$ws.Cells.Item(20, 5).Value = "  +1.48%  "
$ws.Cells.Item(21, 5).Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(22, 4) "4.53"
$ws.Cells.Item(22, 5).Value = "  +2.36%  "
Set-TextValue $ws.Cells.Item(23, 4) "9.48"
$ws.Cells.Item(23, 5).Value = "  +3.06%  "
$ws.Cells.Item(24, 5).Value = "  -3.52%  "
Set-TextValue $ws.Cells.Item(25, 4) "148.07"
$ws.Cells.Item(25, 5).Value = "  +1.70%  "
Set-TextValue $ws.Cells.Item(26, 4) "7.27"
$ws.Cells.Item(26, 5).Value = "  +0.34%  "
Set-TextValue $ws.Cells.Item(27, 4) "16.31"
$ws.Cells.Item(27, 5).Value = "  +2.18%  "
$ws.Cells.Item(28, 5).Value = "  +0.62%  "
$ws.Cells.Item(29, 5).Value = "  +0.12%  "
$ws.Cells.Item(30, 5).Value = "  +0.13%  "
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 4).Value = "1.570.52"
$ws.Cells.Item(32, 5).Value = "  +5.85%  "
$ws.Cells.Item(33, 5).Value = "  +1.63%  "
$ws.Cells.Item(34, 5).Value = "  +2.60%  "
$ws.Cells.Item(35, 5).Value = "  +0.57%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.602"
$ws.Cells.Item(36, 5).Value = "  +3.06%  "
$ws.Cells.Item(37, 5).Value = "  -1.05%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.935"
$ws.Cells.Item(38, 5).Value = "  +4.56%  "
$ws.Cells.Item(39, 5).Value = "  +1.03%  "
$ws.Cells.Item(40, 5).Value = "  +3.67%  "
Set-TextValue $ws.Cells.Item(41, 4) "69.12"
$ws.Cells.Item(41, 5).Value = "  +3.28%  "
$ws.Cells.Item(42, 5).Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(43, 4) "5.59"
$ws.Cells.Item(43, 5).Value = "  -4.76%  "
Set-TextValue $ws.Cells.Item(44, 4) "2.26"
$ws.Cells.Item(44, 5).Value = "  -2.29%  "
$ws.Cells.Item(45, 4).Value = "1.825.61"
$ws.Cells.Item(45, 5).Value = "  +0.46%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.786"
$ws.Cells.Item(46, 5).Value = "  +1.12%  "
Set-TextValue $ws.Cells.Item(47, 4) "90.80"
$ws.Cells.Item(47, 5).Value = "  +0.41%  "
$ws.Cells.Item(48, 5).Value = "  +2.99%  "
$ws.Cells.Item(49, 5).Value = "  +0.77%  "
Set-TextValue $ws.Cells.Item(50, 4) "8.14"
$ws.Cells.Item(50, 5).Value = "  +6.29%  "
$ws.Cells.Item(51, 5).Value = "  +1.80%  "
